$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-03-11"

# Update the header label in I1 to match the new "through" date
$ws.Range("I1").Value = "2022 (through 03-11)"

# Update the data values per the diff
$ws.Range("I3").Value = 141
$ws.Range("I4").Value = 48
$ws.Range("H13").Value = 205
$ws.Range("H14").Value = 1851
$ws.Range("I14").Value = 348
